$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5 from 2023-10-05 to 2023-10-08
$ws.Range("C2:C5").Value = Get-Date -Year 2023 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0
